$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.363
$ws.Range("C2").Value = 0.297
$ws.Range("H2").Value = 2.541
$ws.Range("I2").Value = 2.079
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 74
$ws.Range("N2").Value = 0.2393131875437522
$ws.Range("O2").Value = 0.166830908040689
$ws.Range("Q2").Value = 0.8978683829373846
$ws.Range("R2").Value = 0.8170792379544533

# Row 3
$ws.Range("E3").Value = 7.700000000000001
$ws.Range("F3").Value = 6.3
$ws.Range("H3").Value = 2.541
$ws.Range("I3").Value = 2.079
$ws.Range("K3").Value = 60
$ws.Range("L3").Value = 67
$ws.Range("N3").Value = 0.2393692844097636
$ws.Range("O3").Value = 0.1669656336268933
$ws.Range("Q3").Value = 0.8978681217589644
$ws.Range("R3").Value = 0.8170831300264577

# Row 4
$ws.Range("B4").Value = 0.363
$ws.Range("C4").Value = 0.297
$ws.Range("E4").Value = 7.700000000000001
$ws.Range("F4").Value = 6.3
$ws.Range("H4").Value = 2.795100000000001
$ws.Range("I4").Value = 1.8711
$ws.Range("K4").Value = 52
$ws.Range("L4").Value = 81
$ws.Range("N4").Value = 0.2744746858349662
$ws.Range("O4").Value = 0.1307057700218693
$ws.Range("Q4").Value = 0.9245436742198752
$ws.Range("R4").Value = 0.7577615324283289

$wb.Save()
